$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.08769370117338
$ws.Range("D2").Value = 1.090087753345537
$ws.Range("E2").Value = 1.09002331360952
$ws.Range("F2").Value = 1.10124700516267
$ws.Range("I2").Value = 1.069491623842711
$ws.Range("J2").Value = 1.092537525210105
$ws.Range("K2").Value = 1.092736205283927
$ws.Range("L2").Value = 1.092671930797377
$ws.Range("M2").Value = 1.103867154529329
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.089082650834222
$ws.Range("D3").Value = 1.091243477302077
$ws.Range("E3").Value = 1.091276826003911
$ws.Range("F3").Value = 1.102559642743687
$ws.Range("I3").Value = 1.070024886520648
$ws.Range("J3").Value = 1.093588501226807
$ws.Range("K3").Value = 1.093711178777862
$ws.Range("L3").Value = 1.09374444816225
$ws.Range("M3").Value = 1.105000724982896
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.089980194183306
$ws.Range("D4").Value = 1.091990209655082
$ws.Range("E4").Value = 1.092086550247944
$ws.Range("F4").Value = 1.10340800427791
$ws.Range("I4").Value = 1.070368042182545
$ws.Range("J4").Value = 1.094266829304009
$ws.Range("K4").Value = 1.094340366333619
$ws.Range("L4").Value = 1.094436489356092
$ws.Range("M4").Value = 1.105732658577304
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.090357239665439
$ws.Range("D5").Value = 1.092303876991762
$ws.Range("E5").Value = 1.092426632239522
$ws.Range("F5").Value = 1.103764420108819
$ws.Range("I5").Value = 1.070511852139323
$ws.Range("J5").Value = 1.094551589998475
$ws.Range("K5").Value = 1.094604477091011
$ws.Range("L5").Value = 1.094726961304533
$ws.Range("M5").Value = 1.10603999383335
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.090420530899812
$ws.Range("D6").Value = 1.092356528048977
$ws.Range("E6").Value = 1.092483714519827
$ws.Range("F6").Value = 1.103824250229936
$ws.Range("I6").Value = 1.070535971988043
$ws.Range("J6").Value = 1.094599378744006
$ws.Range("K6").Value = 1.094648799121889
$ws.Range("L6").Value = 1.094775705851037
$ws.Range("M6").Value = 1.106091575227244
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.089985233383666
$ws.Range("D7").Value = 1.091994401905895
$ws.Range("E7").Value = 1.092091095715464
$ws.Range("F7").Value = 1.103412767642632
$ws.Range("I7").Value = 1.070369965553977
$ws.Range("J7").Value = 1.094270635889947
$ws.Range("K7").Value = 1.094343896960258
$ws.Range("L7").Value = 1.094440372466914
$ws.Range("M7").Value = 1.105736766652964
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.088163355009004
$ws.Range("D8").Value = 1.090478565495781
$ws.Range("E8").Value = 1.090447233056691
$ws.Range("F8").Value = 1.10169082764301
$ws.Range("I8").Value = 1.06967223751798
$ws.Range("J8").Value = 1.092893066923424
$ws.Range("K8").Value = 1.093066053436526
$ws.Range("L8").Value = 1.093034799375435
$ws.Range("M8").Value = 1.104250576100674
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.084943533962805
$ws.Range("D9").Value = 1.087798876399323
$ws.Range("E9").Value = 1.087539736383977
$ws.Range("F9").Value = 1.09864865809726
$ws.Range("I9").Value = 1.068428084321367
$ws.Range("J9").Value = 1.090452223236231
$ws.Range("K9").Value = 1.090801249372904
$ws.Range("L9").Value = 1.090542869471421
$ws.Range("M9").Value = 1.10161955600817
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.082790291539148
$ws.Range("D10").Value = 1.086006386852132
$ws.Range("E10").Value = 1.085593840315152
$ws.Range("F10").Value = 1.096614941751934
$ws.Range("I10").Value = 1.067588634806157
$ws.Range("J10").Value = 1.088815723389464
$ws.Range("K10").Value = 1.089282343476879
$ws.Range("L10").Value = 1.088871126390436
$ws.Range("M10").Value = 1.099857087590932
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.081856247571994
$ws.Range("D11").Value = 1.085228731714012
$ws.Range("E11").Value = 1.084749383451912
$ws.Range("F11").Value = 1.095732924199066
$ws.Range("I11").Value = 1.067222734341838
$ws.Range("J11").Value = 1.08810484195476
$ws.Range("K11").Value = 1.088622442699818
$ws.Range("L11").Value = 1.088144699399798
$ws.Range("M11").Value = 1.099091851500841
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.081509043461105
$ws.Range("D12").Value = 1.084939646304387
$ws.Range("E12").Value = 1.084435427519777
$ws.Range("F12").Value = 1.095405086678063
$ws.Range("I12").Value = 1.067086457022163
$ws.Range("J12").Value = 1.087840443207347
$ws.Range("K12").Value = 1.088376990194114
$ws.Range("L12").Value = 1.087874483210818
$ws.Range("M12").Value = 1.098807291440661
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.081583531733422
$ws.Range("D13").Value = 1.08500166657975
$ws.Range("E13").Value = 1.084502785218741
$ws.Range("F13").Value = 1.095475418850089
$ws.Range("I13").Value = 1.067115705562606
$ws.Range("J13").Value = 1.087897173366029
$ws.Range("K13").Value = 1.088429655870355
$ws.Range("L13").Value = 1.087932463187189
$ws.Range("M13").Value = 1.09886834498631
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.081827552886721
$ws.Range("D14").Value = 1.085204840538085
$ws.Range("E14").Value = 1.084723437646993
$ws.Range("F14").Value = 1.095705829506594
$ws.Range("I14").Value = 1.067211477091692
$ws.Range("J14").Value = 1.088082993755154
$ws.Range("K14").Value = 1.088602160400636
$ws.Range("L14").Value = 1.088122371210297
$ws.Range("M14").Value = 1.099068336182128
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.081977867902005
$ws.Range("D15").Value = 1.085329992135976
$ws.Range("E15").Value = 1.084859350674372
$ws.Range("F15").Value = 1.09584776419556
$ws.Range("I15").Value = 1.067270436555801
$ws.Range("J15").Value = 1.088197437859189
$ws.Range("K15").Value = 1.088708401471828
$ws.Range("L15").Value = 1.088239328116084
$ws.Range("M15").Value = 1.099191515148505
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.082852245072519
$ws.Range("D16").Value = 1.086057965330739
$ws.Range("E16").Value = 1.085649844208363
$ws.Range("F16").Value = 1.096673448252684
$ws.Range("I16").Value = 1.067612867311857
$ws.Range("J16").Value = 1.08886285405463
$ws.Range("K16").Value = 1.089326092065071
$ws.Range("L16").Value = 1.088919282731687
$ws.Range("M16").Value = 1.099907829593193
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.08340026530903
$ws.Range("D17").Value = 1.086514199813277
$ws.Range("E17").Value = 1.086145194799401
$ws.Range("F17").Value = 1.097190997781443
$ws.Range("I17").Value = 1.067827016842274
$ws.Range("J17").Value = 1.089279641624419
$ws.Range("K17").Value = 1.089712959486816
$ws.Range("L17").Value = 1.089345113345133
$ws.Range("M17").Value = 1.100356595166452
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.083719755015804
$ws.Range("D18").Value = 1.086780169730763
$ws.Range("E18").Value = 1.086433944296535
$ws.Range("F18").Value = 1.097492740505547
$ws.Range("I18").Value = 1.067951693999116
$ws.Range("J18").Value = 1.089522528464413
$ws.Range("K18").Value = 1.089938400330993
$ws.Range("L18").Value = 1.089593247361876
$ws.Range("M18").Value = 1.10061815263569
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.083828665608351
$ws.Range("D19").Value = 1.086870834335771
$ws.Range("E19").Value = 1.086532370040521
$ws.Range("F19").Value = 1.097595604227403
$ws.Range("I19").Value = 1.06799416631878
$ws.Range("J19").Value = 1.089605309731287
$ws.Range("K19").Value = 1.090015233953675
$ws.Range("L19").Value = 1.089677813129164
$ws.Range("M19").Value = 1.100707303231884
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.083341484651406
$ws.Range("D20").Value = 1.086465265096753
$ws.Range("E20").Value = 1.086092067044369
$ws.Range("F20").Value = 1.097135483639104
$ws.Range("I20").Value = 1.067804064710209
$ws.Range("J20").Value = 1.089234946881591
$ws.Range("K20").Value = 1.089671474276734
$ws.Range("L20").Value = 1.089299451213749
$ws.Range("M20").Value = 1.100308467584741
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0817557019326
$ws.Range("D21").Value = 1.085145017251644
$ws.Range("E21").Value = 1.084658468978575
$ws.Range("F21").Value = 1.095637985346601
$ws.Range("I21").Value = 1.067183284880048
$ws.Range("J21").Value = 1.08802828386827
$ws.Range("K21").Value = 1.088551371429079
$ws.Range("L21").Value = 1.088066458810186
$ws.Range("M21").Value = 1.099009452562809
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.080757157847673
$ws.Range("D22").Value = 1.084313593747886
$ws.Range("E22").Value = 1.083755445196545
$ws.Range("F22").Value = 1.094695190719549
$ws.Range("I22").Value = 1.066790859226767
$ws.Range("J22").Value = 1.087267604020677
$ws.Range("K22").Value = 1.087845172077838
$ws.Range("L22").Value = 1.087288974939132
$ws.Range("M22").Value = 1.098190871706243
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.081286649388852
$ws.Range("D23").Value = 1.08475447494907
$ws.Range("E23").Value = 1.084234314892657
$ws.Range("F23").Value = 1.095195105197542
$ws.Range("I23").Value = 1.06699909314765
$ws.Range("J23").Value = 1.087671046417906
$ws.Range("K23").Value = 1.088219727797427
$ws.Range("L23").Value = 1.087701349238954
$ws.Range("M23").Value = 1.098624993012471
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.083368045598099
$ws.Range("D24").Value = 1.086487377033222
$ws.Range("E24").Value = 1.086116073746176
$ws.Range("F24").Value = 1.097160568507708
$ws.Range("I24").Value = 1.067814436509347
$ws.Range("J24").Value = 1.089255143186169
$ws.Range("K24").Value = 1.089690220313364
$ws.Range("L24").Value = 1.089320084724452
$ws.Range("M24").Value = 1.100330214984585
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.085777089221581
$ws.Range("D25").Value = 1.088492684425885
$ws.Range("E25").Value = 1.088292705006486
$ws.Range("F25").Value = 1.099436097536459
$ws.Range("I25").Value = 1.068751481545062
$ws.Range("J25").Value = 1.091084854174533
$ws.Range("K25").Value = 1.091388330432066
$ws.Range("L25").Value = 1.091188915224357
$ws.Range("M25").Value = 1.102301207566204
